$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 19884
$ws.Cells.Item(2, 2).Value = "Léo Melo"
$ws.Cells.Item(2, 3).Value = "Recursos Humanos"
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(2, 6).Value = 45083
$ws.Cells.Item(2, 7).Value = 7891.64

# Row 3
$ws.Cells.Item(3, 1).Value = 81220
$ws.Cells.Item(3, 2).Value = "João Pedro Sampaio"
$ws.Cells.Item(3, 3).Value = "Vendas"
$ws.Cells.Item(3, 4).Value = "Consulta medica"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 45086
$ws.Cells.Item(3, 7).Value = 5732.94

# Row 4
$ws.Cells.Item(4, 1).Value = 15056
$ws.Cells.Item(4, 2).Value = "Srta. Jade da Paz"
$ws.Cells.Item(4, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(4, 4).Value = "Doenca"
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 45104
$ws.Cells.Item(4, 7).Value = 7663.49

# Row 5
$ws.Cells.Item(5, 1).Value = 43220
$ws.Cells.Item(5, 2).Value = "Enrico Castro"
$ws.Cells.Item(5, 3).Value = "Vendas"
$ws.Cells.Item(5, 4).Value = "Consulta medica"
$ws.Cells.Item(5, 5).Value = 5
$ws.Cells.Item(5, 6).Value = 45090
$ws.Cells.Item(5, 7).Value = 6839.71

# Row 6
$ws.Cells.Item(6, 1).Value = 33882
$ws.Cells.Item(6, 2).Value = "Arthur Gabriel Pacheco"
$ws.Cells.Item(6, 3).Value = "P&D"
$ws.Cells.Item(6, 4).Value = "Doenca"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 45100
$ws.Cells.Item(6, 7).Value = 3301.95

# Row 7
$ws.Cells.Item(7, 1).Value = 11435
$ws.Cells.Item(7, 2).Value = "Vitor Nogueira"
$ws.Cells.Item(7, 3).Value = "P&D"
$ws.Cells.Item(7, 4).Value = "Doenca"
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 45093
$ws.Cells.Item(7, 7).Value = 6240.95

# Row 8
$ws.Cells.Item(8, 1).Value = 24780
$ws.Cells.Item(8, 2).Value = "Heloísa Castro"
$ws.Cells.Item(8, 3).Value = "Engenharia"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 45105
$ws.Cells.Item(8, 7).Value = 3137.22

# Row 9
$ws.Cells.Item(9, 1).Value = 50632
$ws.Cells.Item(9, 2).Value = "Clarice Rios"
$ws.Cells.Item(9, 3).Value = "Atendimento ao Cliente"
$ws.Cells.Item(9, 4).Value = "Consulta medica"
$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = 45096
$ws.Cells.Item(9, 7).Value = 8019.2

# Row 10
$ws.Cells.Item(10, 1).Value = 39489
$ws.Cells.Item(10, 2).Value = "Vinicius Fogaça"
$ws.Cells.Item(10, 5).Value = 5
$ws.Cells.Item(10, 6).Value = 45078
$ws.Cells.Item(10, 7).Value = 7405.24

# Row 11
$ws.Cells.Item(11, 1).Value = 55772
$ws.Cells.Item(11, 2).Value = "Caleb Ramos"
$ws.Cells.Item(11, 3).Value = "Recursos Humanos"
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 45104
$ws.Cells.Item(11, 7).Value = 8689.92
